# minor update on kublet
# Edits slide 4 ("Core Components - worker"): moves the body placeholder up
# slightly, inserts a new bullet after "Monitors the state of the node",
# tweaks the "Pod" bullet wording, and adds a new sub-bullet about Pods
# wrapping containers.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)
$shape = $s.Shapes.Item(1)
$tr = $shape.TextFrame.TextRange

# 1) Shift the placeholder up (504000,1124700) -> (504000,987540) EMU.
#    PowerPoint COM measures Top/Left in points (1 pt = 12700 EMU).
$shape.Top = 987540 / 12700

# 2) Insert a new level-1 bullet "Actually starts containers" right after
#    "Monitors the state of the node" (and before "kube-proxy").
$monitors = $tr.Paragraphs(5, 1)
$null = $monitors.InsertAfter("`rActually starts containers")

# 3) Reword the "Pod" bullet's first run. The earlier insert shifted every
#    paragraph after it down by one, so the detail paragraph is now #12.
$podDetail = $tr.Paragraphs(12, 1)
$firstRun = $podDetail.Runs(1, 1)
$firstRun.Text = "A the smallest, schedulable resource that is managed by "

# 4) Append a new sub-bullet after the Pod detail paragraph describing that
#    Pods wrap around one or more (docker) containers, built up as three
#    separate runs so "docker" stays its own run.
$podDetail = $tr.Paragraphs(12, 1)
$null = $podDetail.InsertAfter("`rPods wrap around one or more (")
$newPara = $tr.Paragraphs(13, 1)
$null = $newPara.InsertAfter("docker")
$newPara = $tr.Paragraphs(13, 1)
$null = $newPara.InsertAfter(") containers")
